$wb = $excel.ActiveWorkbook

# ---- Sheet: Home win ----
$ws = $wb.Worksheets.Item("Home win")
$ws.Cells.Item(2, 1).Value = "30-01-2025 20:00"
$ws.Cells.Item(2, 2).Value = "WORLD"
$ws.Cells.Item(2, 3).Value = "UEFA EUROPA LEAGUE"
$ws.Cells.Item(2, 4).Value = "Dynamo Kyiv - Rīgas FS"
$ws.Cells.Item(2, 5).Value = 73.3
$ws.Cells.Item(2, 6).Value = 1.85
$ws.Cells.Item(3, 1).Value = "30-01-2025 20:00"
$ws.Cells.Item(3, 2).Value = "WORLD"
$ws.Cells.Item(3, 3).Value = "UEFA EUROPA LEAGUE"
$ws.Cells.Item(3, 4).Value = "Maccabi Tel Aviv - FC Porto"
$ws.Cells.Item(3, 5).Value = 70
$ws.Cells.Item(3, 6).Value = 7
$ws.Cells.Item(4, 1).Value = "30-01-2025 20:00"
$ws.Cells.Item(4, 2).Value = "WORLD"
$ws.Cells.Item(4, 3).Value = "UEFA EUROPA LEAGUE"
$ws.Cells.Item(4, 4).Value = "Rangers - Union St. Gilloise"
$ws.Cells.Item(4, 5).Value = 73.3
$ws.Cells.Item(4, 6).Value = 2
$ws.Cells.Item(5, 1).Value = "30-01-2025 20:00"
$ws.Cells.Item(5, 2).Value = "WORLD"
$ws.Cells.Item(5, 3).Value = "UEFA EUROPA LEAGUE"
$ws.Cells.Item(5, 4).Value = "AS Roma - Eintracht Frankfurt"
$ws.Cells.Item(5, 5).Value = 80
$ws.Cells.Item(5, 6).Value = 1.73
$ws.Cells.Item(6, 1).Value = "30-01-2025 23:00"
$ws.Cells.Item(6, 2).Value = "BRAZIL"
$ws.Cells.Item(6, 3).Value = "GAÚCHO - 1"
$ws.Cells.Item(6, 4).Value = "Pelotas - Avenida"
$ws.Cells.Item(6, 5).Value = 73.3
$ws.Cells.Item(6, 6).Value = 2.1
$ws.Cells.Item(7, 1).Value = "30-01-2025 08:30"
$ws.Cells.Item(7, 2).Value = "INDONESIA"
$ws.Cells.Item(7, 3).Value = "LIGA 2"
$ws.Cells.Item(7, 4).Value = "Persikota Tangerang - Sriwijaya FC"
$ws.Cells.Item(7, 5).Value = 80
$ws.Cells.Item(7, 6).Value = 1.91
$ws.Cells.Item(8, 1).Value = "31-01-2025 19:00"
$ws.Cells.Item(8, 2).Value = "NETHERLANDS"
$ws.Cells.Item(8, 3).Value = "EERSTE DIVISIE"
$ws.Cells.Item(8, 4).Value = "Cambuur - Telstar"
$ws.Cells.Item(8, 5).Value = 86.7
$ws.Cells.Item(8, 6).Value = 1.77
$ws.Cells.Item(9, 1).Value = "31-01-2025 18:30"
$ws.Cells.Item(9, 2).Value = "FRANCE"
$ws.Cells.Item(9, 3).Value = "NATIONAL 1"
$ws.Cells.Item(9, 4).Value = "Aubagne - Concarneau"
$ws.Cells.Item(9, 5).Value = 70
$ws.Cells.Item(9, 6).Value = 2.2
$ws.Cells.Item(10, 1).Value = "31-01-2025 18:00"
$ws.Cells.Item(10, 2).Value = "GERMANY"
$ws.Cells.Item(10, 3).Value = "REGIONALLIGA - NORDOST"
$ws.Cells.Item(10, 4).Value = "SV Babelsberg 03 - Hertha Zehlendorf"
$ws.Cells.Item(10, 5).Value = 73.3
$ws.Cells.Item(10, 6).Value = 1.73
$ws.Cells.Item(11, 1).Value = "31-01-2025 18:30"
$ws.Cells.Item(11, 2).Value = "GERMANY"
$ws.Cells.Item(11, 3).Value = "REGIONALLIGA - WEST"
$ws.Cells.Item(11, 4).Value = "MSV Duisburg - Rot-weiss Oberhausen"
$ws.Cells.Item(11, 5).Value = 70
$ws.Cells.Item(11, 6).Value = 1.83
$ws.Cells.Item(12, 1).Value = "31-01-2025 17:00"
$ws.Cells.Item(12, 2).Value = "TURKEY"
$ws.Cells.Item(12, 3).Value = "SÜPER LIG"
$ws.Cells.Item(12, 4).Value = "Eyüpspor - Sivasspor"
$ws.Cells.Item(12, 5).Value = 80
$ws.Cells.Item(12, 6).Value = 1.8

# ---- Sheet: Away Win ----
$ws = $wb.Worksheets.Item("Away Win")
$ws.Cells.Item(2, 1).Value = "31-01-2025 23:30"
$ws.Cells.Item(2, 2).Value = "WORLD"
$ws.Cells.Item(2, 3).Value = "SUDAMERICANO U20"
$ws.Cells.Item(2, 4).Value = "Venezuela U20 - Uruguay U20"
$ws.Cells.Item(2, 5).Value = 86.7
$ws.Cells.Item(2, 6).Value = 1.95
$ws.Rows.Item(3).Delete()

# ---- Sheet: Draw ----
$ws = $wb.Worksheets.Item("Draw")
$ws.Cells.Item(2, 1).Value = "30-01-2025 23:00"
$ws.Cells.Item(2, 2).Value = "BRAZIL"
$ws.Cells.Item(2, 3).Value = "PERNAMBUCANO - 1"
$ws.Cells.Item(2, 4).Value = "Nautico Recife - Maguary PE"
$ws.Cells.Item(2, 5).Value = 60
$ws.Cells.Item(2, 6).Value = 3.5
$ws.Cells.Item(3, 1).Value = "30-01-2025 12:30"
$ws.Cells.Item(3, 2).Value = "IRAN"
$ws.Cells.Item(3, 3).Value = "AZADEGAN LEAGUE"
$ws.Cells.Item(3, 4).Value = "Naft Masjed Soleyman - Shahrdari Noshahr"
$ws.Cells.Item(3, 5).Value = 60
$ws.Cells.Item(3, 6).Value = 2.62
$ws.Cells.Item(4, 1).Value = "31-01-2025 18:00"
$ws.Cells.Item(4, 2).Value = "ROMANIA"
$ws.Cells.Item(4, 3).Value = "LIGA I"
$ws.Cells.Item(4, 4).Value = "Unirea Slobozia - Dinamo Bucuresti"
$ws.Cells.Item(4, 5).Value = 60
$ws.Cells.Item(4, 6).Value = 3.1
$ws.Cells.Item(5, 1).Value = "31-01-2025 15:00"
$ws.Cells.Item(5, 2).Value = "ROMANIA"
$ws.Cells.Item(5, 3).Value = "LIGA I"
$ws.Cells.Item(5, 4).Value = "AFC Hermannstadt - Uta Arad"
$ws.Cells.Item(5, 5).Value = 63.3
$ws.Cells.Item(5, 6).Value = 2.88
$ws.Cells.Item(6, 1).Value = "31-01-2025 00:30"
$ws.Cells.Item(6, 2).Value = "BRAZIL"
$ws.Cells.Item(6, 3).Value = "CARIOCA - 1"
$ws.Cells.Item(6, 4).Value = "Flamengo - Sampaio Corrêa RJ"
$ws.Cells.Item(6, 5).Value = 60
$ws.Cells.Item(6, 6).Value = 7.5
$ws.Cells.Item(7, 1).Value = "31-01-2025 17:00"
$ws.Cells.Item(7, 2).Value = "TURKEY"
$ws.Cells.Item(7, 3).Value = "1. LIG"
$ws.Cells.Item(7, 4).Value = "Yeni Çorumspor - Bandırmaspor"
$ws.Cells.Item(7, 5).Value = 66.7
$ws.Cells.Item(7, 6).Value = 3.1

# ---- Sheet: Btts ----
$ws = $wb.Worksheets.Item("Btts")
$ws.Cells.Item(2, 1).Value = "30-01-2025 20:00"
$ws.Cells.Item(2, 2).Value = "WORLD"
$ws.Cells.Item(2, 3).Value = "UEFA EUROPA LEAGUE"
$ws.Cells.Item(2, 4).Value = "FCSB - Manchester United"
$ws.Cells.Item(2, 5).Value = 76
$ws.Cells.Item(2, 6).Value = 1.8
$ws.Cells.Item(3, 1).Value = "30-01-2025 21:30"
$ws.Cells.Item(3, 2).Value = "BRAZIL"
$ws.Cells.Item(3, 3).Value = "CARIOCA - 1"
$ws.Cells.Item(3, 4).Value = "Boavista SC - Madureira"
$ws.Cells.Item(3, 5).Value = 75.8
$ws.Cells.Item(3, 6).Value = 1.91
$ws.Cells.Item(4, 1).Value = "30-01-2025 18:00"
$ws.Cells.Item(4, 2).Value = "BRAZIL"
$ws.Cells.Item(4, 3).Value = "PERNAMBUCANO - 1"
$ws.Cells.Item(4, 4).Value = "Decisão - Petrolina"
$ws.Cells.Item(4, 5).Value = 80
$ws.Cells.Item(4, 6).Value = 2.2
$ws.Cells.Item(5, 1).Value = "30-01-2025 11:30"
$ws.Cells.Item(5, 2).Value = "IRAQ"
$ws.Cells.Item(5, 3).Value = "IRAQI LEAGUE"
$ws.Cells.Item(5, 4).Value = "Al Karkh - Al Hudod"
$ws.Cells.Item(5, 5).Value = 78.9
$ws.Cells.Item(5, 6).Value = 2.2
$ws.Cells.Item(6, 1).Value = "31-01-2025 19:00"
$ws.Cells.Item(6, 2).Value = "FRANCE"
$ws.Cells.Item(6, 3).Value = "LIGUE 2"
$ws.Cells.Item(6, 4).Value = "Clermont Foot - Ajaccio"
$ws.Cells.Item(6, 5).Value = 76.7
$ws.Cells.Item(6, 6).Value = 2.05
$ws.Cells.Item(7, 1).Value = "31-01-2025 14:30"
$ws.Cells.Item(7, 2).Value = "AZERBAIJAN"
$ws.Cells.Item(7, 3).Value = "PREMYER LIQA"
$ws.Cells.Item(7, 4).Value = "Neftchi Baku - Keshla FC"
$ws.Cells.Item(7, 5).Value = 76.7
$ws.Cells.Item(7, 6).Value = 1.8
$ws.Cells.Item(8, 1).Value = "31-01-2025 00:30"
$ws.Cells.Item(8, 2).Value = "BRAZIL"
$ws.Cells.Item(8, 3).Value = "CARIOCA - 1"
$ws.Cells.Item(8, 4).Value = "Flamengo - Sampaio Corrêa RJ"
$ws.Cells.Item(8, 5).Value = 76
$ws.Cells.Item(8, 6).Value = 3
$ws.Cells.Item(9, 1).Value = "31-01-2025 22:00"
$ws.Cells.Item(9, 2).Value = "CHILE"
$ws.Cells.Item(9, 3).Value = "COPA CHILE"
$ws.Cells.Item(9, 4).Value = "Universidad De Concepcion - Nublense"
$ws.Cells.Item(9, 5).Value = 84
$ws.Cells.Item(9, 6).Value = 1.9
$ws.Cells.Item(10, 1).Value = "31-01-2025 19:30"
$ws.Cells.Item(10, 2).Value = "ITALY"
$ws.Cells.Item(10, 3).Value = "SERIE C - GIRONE C"
$ws.Cells.Item(10, 4).Value = "Giugliano - Foggia"
$ws.Cells.Item(10, 5).Value = 83.3
$ws.Cells.Item(10, 6).Value = 1.73
$ws.Cells.Item(11, 1).Value = "31-01-2025 00:00"
$ws.Cells.Item(11, 2).Value = "JAMAICA"
$ws.Cells.Item(11, 3).Value = "PREMIER LEAGUE"
$ws.Cells.Item(11, 4).Value = "Arnett Gardens - Humble Lions"
$ws.Cells.Item(11, 5).Value = 76.7
$ws.Cells.Item(11, 6).Value = 2.2
$ws.Cells.Item(12, 1).Value = "31-01-2025 19:30"
$ws.Cells.Item(12, 2).Value = "SPAIN"
$ws.Cells.Item(12, 3).Value = "PRIMERA DIVISIÓN RFEF - GROUP 2"
$ws.Cells.Item(12, 4).Value = "Intercity - Marbella"
$ws.Cells.Item(12, 5).Value = 92
$ws.Cells.Item(12, 6).Value = 1.85
$ws.Rows.Item(13).Delete()

# ---- Sheet: Over_Under ----
$ws = $wb.Worksheets.Item("Over_Under")
$ws.Cells.Item(2, 1).Value = "30-01-2025 20:00"
$ws.Cells.Item(2, 2).Value = "WORLD"
$ws.Cells.Item(2, 3).Value = "UEFA EUROPA LEAGUE"
$ws.Cells.Item(2, 4).Value = "AS Roma - Eintracht Frankfurt"
$ws.Cells.Item(2, 5).Value = 80
$ws.Cells.Item(2, 6).Value = 1.8
$ws.Cells.Item(2, 7).Value = 53.3
$ws.Cells.Item(2, 8).Value = 3
$ws.Cells.Item(3, 1).Value = "30-01-2025 20:00"
$ws.Cells.Item(3, 2).Value = "WORLD"
$ws.Cells.Item(3, 3).Value = "UEFA EUROPA LEAGUE"
$ws.Cells.Item(3, 4).Value = "SC Braga - Lazio"
$ws.Cells.Item(3, 5).Value = 86.7
$ws.Cells.Item(3, 6).Value = 1.73
$ws.Cells.Item(3, 7).Value = 33.3
$ws.Cells.Item(3, 8).Value = 2.75
$ws.Cells.Item(4, 1).Value = "31-01-2025 14:30"
$ws.Cells.Item(4, 2).Value = "AZERBAIJAN"
$ws.Cells.Item(4, 3).Value = "PREMYER LIQA"
$ws.Cells.Item(4, 4).Value = "Neftchi Baku - Keshla FC"
$ws.Cells.Item(4, 5).Value = 85
$ws.Cells.Item(4, 6).Value = 1.9
$ws.Cells.Item(4, 7).Value = 40
$ws.Cells.Item(4, 8).Value = 3.25
$ws.Cells.Item(5, 1).Value = "31-01-2025 19:30"
$ws.Cells.Item(5, 2).Value = "ITALY"
$ws.Cells.Item(5, 3).Value = "SERIE C - GIRONE C"
$ws.Cells.Item(5, 4).Value = "Giugliano - Foggia"
$ws.Cells.Item(5, 5).Value = 80
$ws.Cells.Item(5, 6).Value = 1.91
$ws.Cells.Item(5, 7).Value = 45
$ws.Cells.Item(5, 8).Value = 3.5
$ws.Cells.Item(6, 1).Value = "31-01-2025 18:30"
$ws.Cells.Item(6, 2).Value = "SWITZERLAND"
$ws.Cells.Item(6, 3).Value = "CHALLENGE LEAGUE"
$ws.Cells.Item(6, 4).Value = "Stade Lausanne-Ouchy - Étoile Carouge"
$ws.Cells.Item(6, 5).Value = 85
$ws.Cells.Item(6, 6).Value = 1.73
$ws.Cells.Item(6, 7).Value = 65
$ws.Cells.Item(6, 8).Value = 2.88
$ws.Cells.Item(7, 1).Value = "31-01-2025 17:00"
$ws.Cells.Item(7, 2).Value = "TURKEY"
$ws.Cells.Item(7, 3).Value = "SÜPER LIG"
$ws.Cells.Item(7, 4).Value = "Konyaspor - BB Bodrumspor"
$ws.Cells.Item(7, 5).Value = 70
$ws.Cells.Item(7, 6).Value = 2.2
$ws.Cells.Item(7, 7).Value = 65
$ws.Cells.Item(7, 8).Value = 4
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

